# The document has two BTEC/Pearson logo pictures living in the headers
# (header1.xml / header2.xml) and two Pearson logos living in the footers
# (footer1.xml / footer2.xml). The edit renames the picture (wp:docPr /
# pic:cNvPr "name" attribute) on each of those four inline pictures:
#
#   footer "first"   (footer1.xml) : image2.png -> image1.png
#   footer "default" (footer2.xml) : image2.png -> image1.png
#   header "first"   (header1.xml) : image1.jpg -> image2.jpg
#   header "default" (header2.xml) : image1.jpg -> image2.jpg
#
# Word's Sections(1).Headers/Footers collection is 1-based with index 1
# being the primary/"default" story and index 2 being the "first page"
# story, i.e.:
#   Headers(1) -> default header, Headers(2) -> first-page header
#   Footers(1) -> default footer, Footers(2) -> first-page footer

$d = $word.ActiveDocument
$section = $d.Sections(1)

function Rename-LogoInStory($story) {
    if (-not $story.Exists) { return }
    $range = $story.Range
    for ($i = 1; $i -le $range.Paragraphs.Count; $i++) {
        $para = $range.Paragraphs($i)
        $shapes = $para.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                $shape.Name = "image2.jpg"
            } elseif ($shape.AlternativeText -like "*PearsonLogo.png") {
                $shape.Name = "image1.png"
            }
        }
    }
}

# Headers: default (1) and first page (2)
Rename-LogoInStory $section.Headers(1)
Rename-LogoInStory $section.Headers(2)

# Footers: default (1) and first page (2)
Rename-LogoInStory $section.Footers(1)
Rename-LogoInStory $section.Footers(2)
